# Auto-generated Excel COM-interop script to apply cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 45 and 46: Coin/Link swap (FraxShare <-> TheSandbox)
$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'

# Price (D) and Volume(1h) (E) cell updates
# Text-like values are force-written as Text to avoid Excel auto-converting
# numeric-looking strings (e.g. '1.0000') into actual numbers, then the
# cell style is reset back to Normal so no stray style index remains.
$updates = [ordered]@{
    'D2' = '30.190.70'
    'E2' = '  -1.47%  '
    'D3' = '1.857.64'
    'E3' = '  -1.27%  '
    'D4' = '1.0000'
    'E4' = '  -0.07%  '
    'D5' = '232.54'
    'E5' = '  -2.90%  '
    'E6' = '  -0.03%  '
    'D7' = '0.4738'
    'E7' = '  -1.88%  '
    'D8' = '0.2735'
    'E8' = '  -3.55%  '
    'D9' = '0.06411'
    'E9' = '  -1.94%  '
    'D10' = '1.834.09'
    'E10' = '  -7.51%  '
    'D11' = '0.07452'
    'E11' = '  -0.92%  '
    'D12' = '16.26'
    'E12' = '  -2.06%  '
    'D13' = '5.014'
    'E13' = '  -1.86%  '
    'D14' = '85.22'
    'E14' = '  -4.22%  '
    'D15' = '0.6307'
    'E15' = '  -5.44%  '
    'D16' = '30.135.60'
    'E16' = '  -1.51%  '
    'E17' = '  -0.11%  '
    'D18' = '231.25'
    'E18' = '  -0.55%  '
    'D19' = '12.73'
    'E19' = '  -4.71%  '
    'D20' = '0.000007335'
    'E20' = '  -3.77%  '
    'D21' = '2.097.07'
    'E21' = '  -6.37%  '
    'D22' = '1.0000'
    'E22' = '  -0.10%  '
    'D23' = '5.055'
    'E23' = '  -4.63%  '
    'D24' = '0.3976'
    'E24' = '  -4.23%  '
    'D25' = '6.003'
    'E25' = '  -2.95%  '
    'D26' = '9.261'
    'E26' = '  -1.13%  '
    'D27' = '165.25'
    'E27' = '  -1.53%  '
    'D28' = '17.79'
    'E28' = '  -4.64%  '
    'D29' = '1.879'
    'E29' = '  -3.68%  '
    'E30' = '  -2.80%  '
    'D31' = '0.1006'
    'E31' = '  +5.24%  '
    'D32' = '4.140'
    'E32' = '  -5.10%  '
    'D33' = '3.914'
    'E33' = '  -3.38%  '
    'E34' = '  -2.83%  '
    'D35' = '1.140'
    'E35' = '  -6.25%  '
    'E36' = '  -3.52%  '
    'D37' = '1.0000'
    'E37' = '  -0.83%  '
    'E38' = '  -0.36%  '
    'E39' = '  +1.72%  '
    'D40' = '2.636'
    'E40' = '  +0.34%  '
    'D41' = '0.9006'
    'E41' = '  -1.77%  '
    'D42' = '1.954'
    'E42' = '  -6.79%  '
    'D43' = '105.57'
    'E43' = '  -1.01%  '
    'D45' = '0.4095'
    'E45' = '  -4.73%  '
    'D46' = '5.544'
    'E46' = '  -4.87%  '
    'D47' = '7.046'
    'E47' = '  -5.26%  '
    'D48' = '61.04'
    'E48' = '  -5.11%  '
    'D49' = '0.1196'
    'E49' = '  -7.12%  '
    'D50' = '8.736'
    'E50' = '  -2.20%  '
    'D51' = '33.21'
    'E51' = '  -2.07%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = '@'
    $cell.Value = $updates[$ref]
    $cell.Style = 'Normal'
}
